$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 83.375
$ws.Range("I12").Value = 59.57143
$ws.Range("K12").Value = 59.57143
$ws.Range("M12").Value = 110.42857
$ws.Range("H19").Value = 1086.75
$ws.Range("I19").Value = 675
$ws.Range("J19").Value = 1498.5
$ws.Range("K19").Value = 675
$ws.Range("L19").Value = 1498.5
$ws.Range("M19").Value = -500
$ws.Range("N19").Value = -1848.5
$ws.Range("H28").Value = 1776.9
$ws.Range("I28").Value = 2081.0715
$ws.Range("J28").Value = 1067.1666
$ws.Range("K28").Value = 2081.0715
$ws.Range("L28").Value = 1067.1666
$ws.Range("M28").Value = -1596.0715
$ws.Range("N28").Value = -2037.1666
$ws.Range("H43").Value = 2495.75
$ws.Range("I43").Value = 1483
$ws.Range("J43").Value = 2833.3333
$ws.Range("K43").Value = 1483
$ws.Range("L43").Value = 2833.3333
$ws.Range("M43").Value = -1414
$ws.Range("N43").Value = -2971.3333
$ws.Range("H80").Value = 1175
$ws.Range("I80").Value = 1700
$ws.Range("K80").Value = 5100
$ws.Range("M80").Value = -4102
$ws.Range("H83").Value = 1175
$ws.Range("I83").Value = 1700
$ws.Range("K83").Value = 15300
$ws.Range("M83").Value = -10308
$ws.Range("H97").Value = 3253.3076
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 3253.3076
$ws.Range("K97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("M97").Value = 9759.9228
$ws.Range("N97").Value = -10751.9228
$ws.Range("H103").Value = 652.75
$ws.Range("I103").Value = 900
$ws.Range("J103").Value = 570.3333
$ws.Range("K103").Value = 2700
$ws.Range("L103").Value = 1710.9999
$ws.Range("M103").Value = -2114
$ws.Range("N103").Value = -2882.9999
$ws.Range("H116").Value = 15864.571
$ws.Range("I116").Value = 15112.25
$ws.Range("J116").Value = 16867.666
$ws.Range("K116").Value = 15112.25
$ws.Range("L116").Value = 16867.666
$ws.Range("M116").Value = -11670.25
$ws.Range("N116").Value = -23751.666
$ws.Range("H138").Value = 2672.712
$ws.Range("I138").Value = 2192.9285
$ws.Range("K138").Value = 6578.7855
$ws.Range("M138").Value = -1438.7855

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 365
$ws.Range("I5").Value = 365
$ws.Range("K5").Value = 365
$ws.Range("M5").Value = -253
$ws.Range("H16").Value = 16668.334
$ws.Range("J16").Value = 24999.5
$ws.Range("L16").Value = 24999.5
$ws.Range("N16").Value = -25573.5
$ws.Range("H61").Value = 7950.6045
$ws.Range("I61").Value = 6630.657
$ws.Range("J61").Value = 13725.375
$ws.Range("K61").Value = 6630.657
$ws.Range("L61").Value = 13725.375
$ws.Range("M61").Value = -6418.657
$ws.Range("N61").Value = -14149.375
$ws.Range("H82").Value = 85000
$ws.Range("J82").Value = 85000
$ws.Range("L82").Value = 85000
$ws.Range("N82").Value = -85722
$ws.Range("H85").Value = 85000
$ws.Range("J85").Value = 85000
$ws.Range("L85").Value = 85000
$ws.Range("N85").Value = -87496
$ws.Range("H136").Value = 7950.6045
$ws.Range("I136").Value = 6630.657
$ws.Range("J136").Value = 13725.375
$ws.Range("K136").Value = 19891.971
$ws.Range("L136").Value = 41176.125
$ws.Range("M136").Value = -17341.971
$ws.Range("N136").Value = -46276.125

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 365
$ws.Range("I4").Value = 365
$ws.Range("K4").Value = 365
$ws.Range("M4").Value = -250
$ws.Range("H20").Value = 2686.2727
$ws.Range("I20").Value = 1491.6666
$ws.Range("J20").Value = 3513.3076
$ws.Range("K20").Value = 1491.6666
$ws.Range("L20").Value = 3513.3076
$ws.Range("M20").Value = -1244.6666
$ws.Range("N20").Value = -4007.3076

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 2503670.5
$ws.Range("I19").Value = 2503670.5
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 2503670.5
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -2503500.5
$ws.Range("H24").Value = 2503670.5
$ws.Range("I24").Value = 2503670.5
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 2503670.5
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -2503500.5
$ws.Range("H31").Value = 2315.0193
$ws.Range("I31").Value = 1553
$ws.Range("J31").Value = 3275.8262
$ws.Range("K31").Value = 1553
$ws.Range("L31").Value = 3275.8262
$ws.Range("M31").Value = -1258
$ws.Range("N31").Value = -3865.8262
$ws.Range("H34").Value = 2315.0193
$ws.Range("I34").Value = 1553
$ws.Range("J34").Value = 3275.8262
$ws.Range("K34").Value = 1553
$ws.Range("L34").Value = 3275.8262
$ws.Range("M34").Value = -1351
$ws.Range("N34").Value = -3679.8262
$ws.Range("H132").Value = 3145
$ws.Range("I132").Value = 3002.889
$ws.Range("K132").Value = 9008.667000000001
$ws.Range("M132").Value = -6478.667000000001
$ws.Range("H134").Value = 5905.9644
$ws.Range("I134").Value = 5188.72
$ws.Range("J134").Value = 11883
$ws.Range("K134").Value = 15566.16
$ws.Range("L134").Value = 35649
$ws.Range("M134").Value = -13031.16
$ws.Range("N134").Value = -40719

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3126325.2
$ws.Range("I4").Value = 3480712.8
$ws.Range("K4").Value = 10442138.4
$ws.Range("M4").Value = -10442026.4
$ws.Range("H5").Value = 725.4
$ws.Range("I5").Value = 608.2857
$ws.Range("J5").Value = 998.6667
$ws.Range("K5").Value = 1824.8571
$ws.Range("L5").Value = 2996.0001
$ws.Range("M5").Value = -1712.8571
$ws.Range("N5").Value = -3220.0001
$ws.Range("H80").Value = 5356.7144
$ws.Range("J80").Value = 6082.5
$ws.Range("L80").Value = 18247.5
$ws.Range("N80").Value = -20119.5
$ws.Range("H81").Value = 9999.5
$ws.Range("I81").Value = 9999
$ws.Range("K81").Value = 29997
$ws.Range("M81").Value = -28874
$ws.Range("H83").Value = 5356.7144
$ws.Range("J83").Value = 6082.5
$ws.Range("L83").Value = 54742.5
$ws.Range("N83").Value = -64102.5
$ws.Range("H84").Value = 9999.5
$ws.Range("I84").Value = 9999
$ws.Range("K84").Value = 89991
$ws.Range("M84").Value = -84375
$ws.Range("H125").Value = 14333.333
$ws.Range("J125").Value = 15333.333
$ws.Range("L125").Value = 45999.999
$ws.Range("N125").Value = -55839.999
$ws.Range("H135").Value = 725.4
$ws.Range("I135").Value = 608.2857
$ws.Range("J135").Value = 998.6667
$ws.Range("K135").Value = 5474.571300000001
$ws.Range("L135").Value = 8988.0003
$ws.Range("M135").Value = -2939.571300000001
$ws.Range("N135").Value = -14058.0003
$ws.Range("H140").Value = 5008117.5
$ws.Range("I140").Value = 13160159
$ws.Range("K140").Value = 39480477
$ws.Range("M140").Value = -39475297

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 18648.5
$ws.Range("I80").Value = 1943
$ws.Range("J80").Value = 27001.25
$ws.Range("K80").Value = 1943
$ws.Range("L80").Value = 27001.25
$ws.Range("M80").Value = -945
$ws.Range("N80").Value = -28997.25
$ws.Range("H83").Value = 18648.5
$ws.Range("I83").Value = 1943
$ws.Range("J83").Value = 27001.25
$ws.Range("K83").Value = 9715
$ws.Range("L83").Value = 135006.25
$ws.Range("M83").Value = -4723
$ws.Range("N83").Value = -144990.25
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("N103").Value = 0
$ws.Range("H126").Value = 6272.75
$ws.Range("I126").Value = 4028.8333
$ws.Range("J126").Value = 8516.666999999999
$ws.Range("K126").Value = 12086.4999
$ws.Range("L126").Value = 25550.001
$ws.Range("M126").Value = -9616.499899999999
$ws.Range("N126").Value = -30490.001

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2750
$ws.Range("J7").Value = 2750
$ws.Range("L7").Value = 2750
$ws.Range("N7").Value = -2974
$ws.Range("H18").Value = 55600
$ws.Range("J18").Value = 48333.332
$ws.Range("L18").Value = 48333.332
$ws.Range("N18").Value = -48677.332
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H61").Value = 2986
$ws.Range("I61").Value = 3360.4443
$ws.Range("J61").Value = 1862.6666
$ws.Range("K61").Value = 3360.4443
$ws.Range("L61").Value = 1862.6666
$ws.Range("M61").Value = -3158.4443
$ws.Range("N61").Value = -2266.6666
$ws.Range("H82").Value = 6843.5557
$ws.Range("I82").Value = 1934
$ws.Range("J82").Value = 9298.333000000001
$ws.Range("K82").Value = 1934
$ws.Range("L82").Value = 9298.333000000001
$ws.Range("M82").Value = -1573
$ws.Range("N82").Value = -10020.333
$ws.Range("H85").Value = 6843.5557
$ws.Range("I85").Value = 1934
$ws.Range("J85").Value = 9298.333000000001
$ws.Range("K85").Value = 1934
$ws.Range("L85").Value = 9298.333000000001
$ws.Range("M85").Value = -686
$ws.Range("N85").Value = -11794.333
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H113").Value = 2986
$ws.Range("I113").Value = 3360.4443
$ws.Range("J113").Value = 1862.6666
$ws.Range("K113").Value = 3360.4443
$ws.Range("L113").Value = 1862.6666
$ws.Range("M113").Value = -1190.4443
$ws.Range("N113").Value = -6202.6666
$ws.Range("H126").Value = 2750
$ws.Range("J126").Value = 2750
$ws.Range("L126").Value = 8250
$ws.Range("N126").Value = -13190
$ws.Range("H132").Value = 2178.3408
$ws.Range("J132").Value = 3581.125
$ws.Range("L132").Value = 10743.375
$ws.Range("N132").Value = -15803.375
$ws.Range("H136").Value = 1855.8511
$ws.Range("I136").Value = 1052.3125
$ws.Range("K136").Value = 3156.9375
$ws.Range("M136").Value = -606.9375

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4874.3384
$ws.Range("I136").Value = 5009.1187
$ws.Range("K136").Value = 15027.3561
$ws.Range("M136").Value = -12477.3561
